$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Copy the "Premium" highlight format (yellow fill, style used on E2) onto
#    the rows that become premium ("Yes") in the new data (E5 and E7), before
#    any values are touched. E2 keeps the same formatting in the new data too.
# ---------------------------------------------------------------------------
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Make sure the ID column keeps being stored as text (it was inline text
#    in the original file, not a number) before writing the new numeric-
#    looking IDs into column A.
# ---------------------------------------------------------------------------
$ws.Range("A2:A16").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 3) Overwrite rows 2-16 with the refreshed scrape data.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "1327967"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1327967"
$ws.Range("C2").Value = "ACE Program | Russian Financial Analyst"
$ws.Range("D2").Value = "Thane, Maharashtra, India"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "Tata Consultancy Services Ltd."

$ws.Range("A3").Value = "1327988"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1327988"
$ws.Range("C3").Value = "[Impact Brazil] - Sales Intern"
$ws.Range("D3").Value = "Curitiba, PR, Brasil"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "The Coffee"

$ws.Range("A4").Value = "1327981"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1327981"
$ws.Range("C4").Value = "Medical Digital and Administrative Project Analyst"
$ws.Range("D4").Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "4 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "Merck"

$ws.Range("A5").Value = "1327977"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1327977"
$ws.Range("C5").Value = "Taste Hungary |  Customer Services Associate (Portuguese Speaker)"
$ws.Range("D5").Value = "Budapeste, Hungria"
$ws.Range("E5").Value = "Yes"
$ws.Range("F5").Value = "5 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "EATON"

$ws.Range("A6").Value = "1327970"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1327970"
$ws.Range("C6").Value = "Business Development Intern"
$ws.Range("D6").Value = "Malabe, Sri Lanka"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "0 applicants"
$ws.Range("G6").Value = "3 - 6 Months"
$ws.Range("H6").Value = "Ribelz Integrated Pvt Ltd"

$ws.Range("A7").Value = "1327965"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1327965"
$ws.Range("C7").Value = "ACE Program | German Financial Analyst"
$ws.Range("D7").Value = "Thane, Maharashtra, India"
$ws.Range("E7").Value = "Yes"
$ws.Range("F7").Value = "0 applicants"
$ws.Range("G7").Value = "6 - 18 Months"
$ws.Range("H7").Value = "Tata Consultancy Services Ltd."

$ws.Range("A8").Value = "1327882"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1327882"
$ws.Range("C8").Value = "Web developer"
$ws.Range("D8").Value = "El-Mahalla El-Kubra, Al Mahalah Al Kubra (Part 2), El Mahalla El Kubra, Gharbia Governorate, Egypt"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "5 applicants"
$ws.Range("G8").Value = "9 - 12 Weeks"
$ws.Range("H8").Value = "Positive Kids academy"

$ws.Range("A9").Value = "1327128"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1327128"
$ws.Range("C9").Value = "Digital Marketing"
$ws.Range("D9").Value = "Καλλιθέα 630 77, Ελλάδα"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "36 applicants"
$ws.Range("G9").Value = "9 - 12 Weeks"
$ws.Range("H9").Value = "Respirotours"

$ws.Range("A10").Value = "1327043"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1327043"
$ws.Range("C10").Value = "Web Developer"
$ws.Range("D10").Value = "Sousse, Tunisia"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "19 applicants"
$ws.Range("G10").Value = "9 - 12 Weeks"
$ws.Range("H10").Value = "Progress Professional Center"

$ws.Range("A11").Value = "1326995"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1326995"
$ws.Range("C11").Value = "Sales Specialist"
$ws.Range("D11").Value = "10th of Ramadan City, Al-Sharqia Governorate, Egypt"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "7 applicants"
$ws.Range("G11").Value = "9 - 12 Weeks"
$ws.Range("H11").Value = "ABD Eldaem Road Transportation Company"

$ws.Range("A12").Value = "1326990"
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1326990"
$ws.Range("C12").Value = "Artificial intelligence"
$ws.Range("D12").Value = "10th of Ramadan City, Al-Sharqia Governorate, Egypt"
$ws.Range("E12").Value = "No"
$ws.Range("F12").Value = "32 applicants"
$ws.Range("G12").Value = "9 - 12 Weeks"
$ws.Range("H12").Value = "ABD Eldaem Road Transportation Company"

$ws.Range("A13").Value = "1326658"
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1326658"
$ws.Range("C13").Value = "Business Development Intern (Japanese Speaking Individuals Only)"
$ws.Range("D13").Value = "Malabe, Sri Lanka"
$ws.Range("E13").Value = "No"
$ws.Range("F13").Value = "19 applicants"
$ws.Range("G13").Value = "3 - 6 Months"
$ws.Range("H13").Value = "Creative Technology Solutions (Private) Limited"

$ws.Range("A14").Value = "1317292"
$ws.Range("B14").Value = "https://aiesec.org/opportunity/global-talent/1317292"
$ws.Range("C14").Value = "[Impact Florianópolis]- Social Media"
$ws.Range("D14").Value = "São Miguel do Oeste, SC, 89900-000, Brasil"
$ws.Range("E14").Value = "No"
$ws.Range("F14").Value = "76 applicants"
$ws.Range("G14").Value = "9 - 12 Weeks"
$ws.Range("H14").Value = "KNN Idiomas"

$ws.Range("A15").Value = "1302356"
$ws.Range("B15").Value = "https://aiesec.org/opportunity/global-talent/1302356"
$ws.Range("C15").Value = "Language Specialist - French"
$ws.Range("D15").Value = "Colombo, Sri Lanka"
$ws.Range("E15").Value = "No"
$ws.Range("F15").Value = "28 applicants"
$ws.Range("G15").Value = "3 - 6 Months"
$ws.Range("H15").Value = "Aitken Spence Travels (Pvt) Ltd"

$ws.Range("A16").Value = "1289377"
$ws.Range("B16").Value = "https://aiesec.org/opportunity/global-talent/1289377"
$ws.Range("C16").Value = "Medical Advisor (Italian Speaker)"
$ws.Range("D16").Value = "İstanbul, Türkiye"
$ws.Range("E16").Value = "No"
$ws.Range("F16").Value = "38 applicants"
$ws.Range("G16").Value = "6 - 18 Months"
$ws.Range("H16").Value = "International Plus"

# ---------------------------------------------------------------------------
# 4) Row 17 (old "Medical Advisor (Spanish Speaker)" listing) no longer
#    exists in the refreshed scrape - delete it entirely so the sheet shrinks
#    from 17 to 16 data rows, and the dimension updates automatically.
# ---------------------------------------------------------------------------
$ws.Rows.Item(17).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 5) Column width adjustments (C, D, F, H). Excel's ColumnWidth property adds
#    a constant ~0.8333 character padding once the value round-trips through
#    the stored OOXML width, so we subtract that back out to land exactly on
#    the target widths of 68 / 101 / 16 / 50.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 67.16666666666667
$ws.Columns.Item(4).ColumnWidth = 100.16666666666667
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 49.166666666666664
